# "Exception handling and Encoder decoder Implementation"
#
# The DATA sheet's password column (E) used to store plaintext passwords
# ("admin123" / "admin1234"). They are replaced here with a single
# Base64-encoded value ("admin123" -> "YWRtaW4xMjM=") for every data row,
# entered with a leading apostrophe so Excel stores/treats it as literal
# text (quote-prefixed) rather than re-interpreting it. Once the old
# plaintext strings are no longer referenced anywhere, the shared-string
# table naturally loses them when the workbook is saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DATA")

# Replace the plaintext passwords with the Base64-encoded credential.
# Leading "'" forces text/quote-prefix storage for the value.
$ws.Range("E2").Value = "'YWRtaW4xMjM="
$ws.Range("E3").Value = "'YWRtaW4xMjM="
$ws.Range("E4").Value = "'YWRtaW4xMjM="

# Widen column E so the longer encoded value is fully visible, then select
# it (mirrors clicking the column header), which also makes DATA the
# active sheet/tab.
$ws.Columns("E").ColumnWidth = 14.25
[void]$ws.Columns("E").Select()
[void]$ws.Activate()
